$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the F:V data between row 55 and row 56 (A:E stay untouched) ---
$row55 = $ws.Range("F55:V55").Value2
$row56 = $ws.Range("F56:V56").Value2

$ws.Range("F55:V55").Value2 = $row56
$ws.Range("F56:V56").Value2 = $row55

# --- Append the new match as row 76 ---
$ws.Range("A76").Value2 = 75
$ws.Range("B76").Value2 = "armenia"
$ws.Range("C76").Value2 = "premier-league"
$ws.Range("D76").Value2 = "2023-2024"
$ws.Range("E76").Value2 = 45238.66666666666
$ws.Range("F76").Value2 = "Urartu"
$ws.Range("G76").Value2 = 1
$ws.Range("H76").Value2 = "Alashkert"
$ws.Range("I76").Value2 = 0
$ws.Range("J76").Value2 = 2.07
$ws.Range("K76").Value2 = "07/11/2023 04:12"
$ws.Range("L76").Value2 = 2.26
$ws.Range("M76").Value2 = "08/11/2023 15:58"
$ws.Range("N76").Value2 = 3.38
$ws.Range("O76").Value2 = "07/11/2023 04:12"
$ws.Range("P76").Value2 = 3.53
$ws.Range("Q76").Value2 = "08/11/2023 15:58"
$ws.Range("R76").Value2 = 3.1
$ws.Range("S76").Value2 = "07/11/2023 04:12"
$ws.Range("T76").Value2 = 3.01
$ws.Range("U76").Value2 = "08/11/2023 15:58"
$ws.Range("V76").Value2 = "https://www.betexplorer.com/football/armenia/premier-league/urartu-alashkert/x0bprFkt/"

# Match formatting of the preceding data row (bold/border index style for A, date style for E)
$ws.Range("A75").Copy()
$ws.Range("A76").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E75").Copy()
$ws.Range("E76").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
